$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value2 = [double]"0.08315501362085342"
$ws.Cells.Item(2, 2).Value2 = [double]"0.9773325324058533"
$ws.Cells.Item(2, 3).Value2 = [double]"0.008656255900859833"
$ws.Cells.Item(2, 4).Value2 = [double]"0.9987184405326843"
$ws.Cells.Item(3, 1).Value2 = [double]"0.0126164061948657"
$ws.Cells.Item(3, 2).Value2 = [double]"0.9980259537696838"
$ws.Cells.Item(3, 3).Value2 = [double]"0.005805364344269037"
$ws.Cells.Item(3, 4).Value2 = [double]"0.9990388751029968"
$ws.Cells.Item(4, 1).Value2 = [double]"0.008325016126036644"
$ws.Cells.Item(4, 2).Value2 = [double]"0.9983209371566772"
$ws.Cells.Item(4, 3).Value2 = [double]"0.002757597947493196"
$ws.Cells.Item(4, 4).Value2 = [double]"0.9993592500686646"
$ws.Cells.Item(5, 1).Value2 = [double]"0.003220483893528581"
$ws.Cells.Item(5, 2).Value2 = [double]"0.9992512464523315"
$ws.Cells.Item(5, 3).Value2 = [double]"0.0005781799554824829"
$ws.Cells.Item(5, 4).Value2 = [double]"0.9997597336769104"
$ws.Cells.Item(6, 1).Value2 = [double]"0.003856554394587874"
$ws.Cells.Item(6, 2).Value2 = [double]"0.9992738962173462"
$ws.Cells.Item(6, 3).Value2 = [double]"0.0002524404844734818"
$ws.Cells.Item(6, 4).Value2 = [double]"0.9998397827148438"
$ws.Cells.Item(7, 1).Value2 = [double]"0.001940985908731818"
$ws.Cells.Item(7, 2).Value2 = [double]"0.9996142387390137"
$ws.Cells.Item(7, 3).Value2 = [double]"0.0007461768691428006"
$ws.Cells.Item(7, 4).Value2 = [double]"0.9996796250343323"
$ws.Cells.Item(8, 1).Value2 = [double]"0.001462956657633185"
$ws.Cells.Item(8, 2).Value2 = [double]"0.9997504353523254"
$ws.Cells.Item(8, 3).Value2 = [double]"0.0003442996821831912"
$ws.Cells.Item(8, 4).Value2 = [double]"0.9998397827148438"
$ws.Cells.Item(9, 1).Value2 = [double]"0.001754448283463717"
$ws.Cells.Item(9, 2).Value2 = [double]"0.9996596574783325"
$ws.Cells.Item(9, 3).Value2 = [double]"0.0006934424745850265"
$ws.Cells.Item(9, 4).Value2 = [double]"0.9996796250343323"
$ws.Cells.Item(10, 1).Value2 = [double]"0.0006591174169443548"
$ws.Cells.Item(10, 2).Value2 = [double]"0.9997730851173401"
$ws.Cells.Item(10, 3).Value2 = [double]"0.0005000810488127172"
$ws.Cells.Item(10, 4).Value2 = [double]"0.9998397827148438"
$ws.Cells.Item(11, 1).Value2 = [double]"0.00131720129866153"
$ws.Cells.Item(11, 2).Value2 = [double]"0.9996142387390137"
$ws.Cells.Item(11, 3).Value2 = [double]"2.675762698345352E-05"
$ws.Cells.Item(11, 4).Value2 = [double]"1"
$ws.Cells.Item(12, 1).Value2 = [double]"0.0006066200439818203"
$ws.Cells.Item(12, 2).Value2 = [double]"0.999863862991333"
$ws.Cells.Item(12, 3).Value2 = [double]"0.0003097179578617215"
$ws.Cells.Item(12, 4).Value2 = [double]"0.9998397827148438"
$ws.Cells.Item(13, 1).Value2 = [double]"0.0006012291414663196"
$ws.Cells.Item(13, 2).Value2 = [double]"0.9997730851173401"
$ws.Cells.Item(13, 3).Value2 = [double]"0.0001013177388813347"
$ws.Cells.Item(13, 4).Value2 = [double]"1"
$ws.Cells.Item(14, 1).Value2 = [double]"0.0008342261426150799"
$ws.Cells.Item(14, 2).Value2 = [double]"0.9998865723609924"
$ws.Cells.Item(14, 3).Value2 = [double]"7.514948083553463E-05"
$ws.Cells.Item(14, 4).Value2 = [double]"1"
$ws.Cells.Item(15, 1).Value2 = [double]"5.136373147252016E-05"
$ws.Cells.Item(15, 2).Value2 = [double]"1"
$ws.Cells.Item(15, 3).Value2 = [double]"6.173732344905147E-06"
$ws.Cells.Item(15, 4).Value2 = [double]"1"
$ws.Cells.Item(16, 1).Value2 = [double]"0.0007021779892966151"
$ws.Cells.Item(16, 2).Value2 = [double]"0.9997957944869995"
$ws.Cells.Item(16, 3).Value2 = [double]"1.363915203000943E-06"
$ws.Cells.Item(16, 4).Value2 = [double]"1"
$ws.Cells.Item(17, 1).Value2 = [double]"0.000599684368353337"
$ws.Cells.Item(17, 2).Value2 = [double]"0.9998865723609924"
$ws.Cells.Item(17, 3).Value2 = [double]"6.828194045738201E-07"
$ws.Cells.Item(17, 4).Value2 = [double]"1"
$ws.Cells.Item(18, 1).Value2 = [double]"4.516240733209997E-05"
$ws.Cells.Item(18, 2).Value2 = [double]"0.9999772906303406"
$ws.Cells.Item(18, 3).Value2 = [double]"5.34313983280299E-07"
$ws.Cells.Item(18, 4).Value2 = [double]"1"
$ws.Cells.Item(19, 1).Value2 = [double]"0.0003091931575909257"
$ws.Cells.Item(19, 2).Value2 = [double]"0.999863862991333"
$ws.Cells.Item(19, 3).Value2 = [double]"0.000134128742502071"
$ws.Cells.Item(19, 4).Value2 = [double]"0.9999198913574219"
$ws.Cells.Item(20, 1).Value2 = [double]"0.001270797802135348"
$ws.Cells.Item(20, 2).Value2 = [double]"0.9996596574783325"
$ws.Cells.Item(20, 3).Value2 = [double]"2.657196273503359E-05"
$ws.Cells.Item(20, 4).Value2 = [double]"1"
$ws.Cells.Item(21, 1).Value2 = [double]"0.0003688965225592256"
$ws.Cells.Item(21, 2).Value2 = [double]"0.9999319314956665"
$ws.Cells.Item(21, 3).Value2 = [double]"2.224756059376887E-07"
$ws.Cells.Item(21, 4).Value2 = [double]"1"
$ws.Cells.Item(22, 1).Value2 = [double]"0.0006585626979358494"
$ws.Cells.Item(22, 2).Value2 = [double]"0.9998865723609924"
$ws.Cells.Item(22, 3).Value2 = [double]"9.305878734267026E-08"
$ws.Cells.Item(22, 4).Value2 = [double]"1"
$ws.Cells.Item(23, 1).Value2 = [double]"0.0002061673876596615"
$ws.Cells.Item(23, 2).Value2 = [double]"0.9999546408653259"
$ws.Cells.Item(23, 3).Value2 = [double]"1.54807196395268E-07"
$ws.Cells.Item(23, 4).Value2 = [double]"1"
$ws.Cells.Item(24, 1).Value2 = [double]"0.0007675917586311698"
$ws.Cells.Item(24, 2).Value2 = [double]"0.9998185038566589"
$ws.Cells.Item(24, 3).Value2 = [double]"5.78155976427297E-08"
$ws.Cells.Item(24, 4).Value2 = [double]"1"
$ws.Cells.Item(25, 1).Value2 = [double]"0.0006697191274724901"
$ws.Cells.Item(25, 2).Value2 = [double]"0.9999546408653259"
$ws.Cells.Item(25, 3).Value2 = [double]"1.269851622964779E-06"
$ws.Cells.Item(25, 4).Value2 = [double]"1"
$ws.Cells.Item(26, 1).Value2 = [double]"0.0003406107716728002"
$ws.Cells.Item(26, 2).Value2 = [double]"0.9998865723609924"
$ws.Cells.Item(26, 3).Value2 = [double]"4.68504737227704E-08"
$ws.Cells.Item(26, 4).Value2 = [double]"1"
$ws.Cells.Item(27, 1).Value2 = [double]"0.000675029878038913"
$ws.Cells.Item(27, 2).Value2 = [double]"0.9998865723609924"
$ws.Cells.Item(27, 3).Value2 = [double]"1.498220001394657E-07"
$ws.Cells.Item(27, 4).Value2 = [double]"1"
$ws.Cells.Item(28, 1).Value2 = [double]"4.625549263437279E-05"
$ws.Cells.Item(28, 2).Value2 = [double]"1"
$ws.Cells.Item(28, 3).Value2 = [double]"1.052818205948824E-07"
$ws.Cells.Item(28, 4).Value2 = [double]"1"
$ws.Cells.Item(29, 1).Value2 = [double]"0.0004655005177482963"
$ws.Cells.Item(29, 2).Value2 = [double]"0.9998411536216736"
$ws.Cells.Item(29, 3).Value2 = [double]"3.590131882447167E-06"
$ws.Cells.Item(29, 4).Value2 = [double]"1"
$ws.Cells.Item(30, 1).Value2 = [double]"0.0001992036413867027"
$ws.Cells.Item(30, 2).Value2 = [double]"0.9999546408653259"
$ws.Cells.Item(30, 3).Value2 = [double]"3.423273483349476E-06"
$ws.Cells.Item(30, 4).Value2 = [double]"1"
$ws.Cells.Item(31, 1).Value2 = [double]"4.922739026369527E-05"
$ws.Cells.Item(31, 2).Value2 = [double]"1"
$ws.Cells.Item(31, 3).Value2 = [double]"1.676548677664869E-08"
$ws.Cells.Item(31, 4).Value2 = [double]"1"
$ws.Cells.Item(32, 1).Value2 = [double]"0.0001970739831449464"
$ws.Cells.Item(32, 2).Value2 = [double]"0.9999092221260071"
$ws.Cells.Item(32, 3).Value2 = [double]"0.002547807991504669"
$ws.Cells.Item(32, 4).Value2 = [double]"0.9998397827148438"
$ws.Cells.Item(33, 1).Value2 = [double]"0.0007342449971474707"
$ws.Cells.Item(33, 2).Value2 = [double]"0.999863862991333"
$ws.Cells.Item(33, 3).Value2 = [double]"3.023961880899151E-06"
$ws.Cells.Item(33, 4).Value2 = [double]"1"
$ws.Cells.Item(34, 1).Value2 = [double]"1.857917959569022E-05"
$ws.Cells.Item(34, 2).Value2 = [double]"1"
$ws.Cells.Item(34, 3).Value2 = [double]"1.131363546846842E-06"
$ws.Cells.Item(34, 4).Value2 = [double]"1"
$ws.Cells.Item(35, 1).Value2 = [double]"0.0004506352997850627"
$ws.Cells.Item(35, 2).Value2 = [double]"0.9999092221260071"
$ws.Cells.Item(35, 3).Value2 = [double]"1.002831950813743E-07"
$ws.Cells.Item(35, 4).Value2 = [double]"1"
$ws.Cells.Item(36, 1).Value2 = [double]"0.0003837404365185648"
$ws.Cells.Item(36, 2).Value2 = [double]"0.9998411536216736"
$ws.Cells.Item(36, 3).Value2 = [double]"6.110747907683844E-09"
$ws.Cells.Item(36, 4).Value2 = [double]"1"
$ws.Cells.Item(37, 1).Value2 = [double]"0.0003577698662411422"
$ws.Cells.Item(37, 2).Value2 = [double]"0.9999319314956665"
$ws.Cells.Item(37, 3).Value2 = [double]"3.34186817108062E-10"
$ws.Cells.Item(37, 4).Value2 = [double]"1"
$ws.Cells.Item(38, 1).Value2 = [double]"0.0003000157012138516"
$ws.Cells.Item(38, 2).Value2 = [double]"0.9999092221260071"
$ws.Cells.Item(38, 3).Value2 = [double]"1.422678641560537E-09"
$ws.Cells.Item(38, 4).Value2 = [double]"1"
$ws.Cells.Item(39, 1).Value2 = [double]"4.64539771201089E-05"
$ws.Cells.Item(39, 2).Value2 = [double]"0.9999772906303406"
$ws.Cells.Item(39, 3).Value2 = [double]"4.201206604825813E-10"
$ws.Cells.Item(39, 4).Value2 = [double]"1"
$ws.Cells.Item(40, 1).Value2 = [double]"0.0003385647432878613"
$ws.Cells.Item(40, 2).Value2 = [double]"0.9999319314956665"
$ws.Cells.Item(40, 3).Value2 = [double]"0.000100694815046154"
$ws.Cells.Item(40, 4).Value2 = [double]"1"
$ws.Cells.Item(41, 1).Value2 = [double]"7.001096673775464E-05"
$ws.Cells.Item(41, 2).Value2 = [double]"0.9999772906303406"
$ws.Cells.Item(41, 3).Value2 = [double]"1.102543478737061E-06"
$ws.Cells.Item(41, 4).Value2 = [double]"1"
$ws.Cells.Item(42, 1).Value2 = [double]"0.0005310841370373964"
$ws.Cells.Item(42, 2).Value2 = [double]"0.999863862991333"
$ws.Cells.Item(42, 3).Value2 = [double]"0.0002739960036706179"
$ws.Cells.Item(42, 4).Value2 = [double]"0.9998397827148438"
$ws.Cells.Item(43, 1).Value2 = [double]"0.001024087541736662"
$ws.Cells.Item(43, 2).Value2 = [double]"0.999863862991333"
$ws.Cells.Item(43, 3).Value2 = [double]"0.0002763093798421323"
$ws.Cells.Item(43, 4).Value2 = [double]"0.9998397827148438"
$ws.Cells.Item(44, 1).Value2 = [double]"1.942887320183218E-05"
$ws.Cells.Item(44, 2).Value2 = [double]"1"
$ws.Cells.Item(44, 3).Value2 = [double]"0.0001300514850299805"
$ws.Cells.Item(44, 4).Value2 = [double]"0.9998397827148438"
$ws.Cells.Item(45, 1).Value2 = [double]"0.0001620481634745374"
$ws.Cells.Item(45, 2).Value2 = [double]"0.9999319314956665"
$ws.Cells.Item(45, 3).Value2 = [double]"0.0002705890219658613"
$ws.Cells.Item(45, 4).Value2 = [double]"0.9998397827148438"
$ws.Cells.Item(46, 1).Value2 = [double]"2.624984699650668E-05"
$ws.Cells.Item(46, 2).Value2 = [double]"1"
$ws.Cells.Item(46, 3).Value2 = [double]"2.349628994124942E-05"
$ws.Cells.Item(46, 4).Value2 = [double]"1"
$ws.Cells.Item(47, 1).Value2 = [double]"6.959089660085738E-06"
$ws.Cells.Item(47, 2).Value2 = [double]"1"
$ws.Cells.Item(47, 3).Value2 = [double]"0.0002268574462505057"
$ws.Cells.Item(47, 4).Value2 = [double]"0.9998397827148438"
$ws.Cells.Item(48, 1).Value2 = [double]"0.0008469157037325203"
$ws.Cells.Item(48, 2).Value2 = [double]"0.9998185038566589"
$ws.Cells.Item(48, 3).Value2 = [double]"7.208343788533966E-08"
$ws.Cells.Item(48, 4).Value2 = [double]"1"
$ws.Cells.Item(49, 1).Value2 = [double]"0.0004613042110577226"
$ws.Cells.Item(49, 2).Value2 = [double]"0.9999546408653259"
$ws.Cells.Item(49, 3).Value2 = [double]"1.574762364953131E-07"
$ws.Cells.Item(49, 4).Value2 = [double]"1"
$ws.Cells.Item(50, 1).Value2 = [double]"9.758869418874383E-05"
$ws.Cells.Item(50, 2).Value2 = [double]"0.9999772906303406"
$ws.Cells.Item(50, 3).Value2 = [double]"2.080440572171938E-05"
$ws.Cells.Item(50, 4).Value2 = [double]"1"
$ws.Cells.Item(51, 1).Value2 = [double]"4.725187682197429E-05"
$ws.Cells.Item(51, 2).Value2 = [double]"0.9999772906303406"
$ws.Cells.Item(51, 3).Value2 = [double]"2.278741931149852E-06"
$ws.Cells.Item(51, 4).Value2 = [double]"1"
